$wb = $excel.ActiveWorkbook

# --- New "2022-Q1" sheet ----------------------------------------------------
# Duplicate the "2021-Q4" sheet (same headers/column layout/styles) and place
# the copy immediately before the last sheet ("总计"), which is where the new
# quarter belongs in the tab order.
$src = $wb.Worksheets.Item("2021-Q4")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($lastSheet)
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

# The template only has 6 data rows (rows 2-7); 2022-Q1 needs 12 (rows 2-13).
# Extend the table by copying the formatting of the last template row down.
$new.Range("A7:H7").Copy()
$new.Range("A8:H13").PasteSpecial(-4122)

# Columns B:G hold fund codes / figures kept as plain text in this workbook
# (e.g. "001371", "40.15" - leading zeros / fixed decimals must survive), so
# force text formatting before writing, then clear the formatting back off
# again so the cells end up as plain, unstyled text (matching the rest of
# the table) rather than keeping a lingering "@" number format.
$new.Range("B2:G13").NumberFormat = "@"

$new.Range("A2").Value = 0
$new.Range("B2").Value = '001371'
$new.Range("C2").Value = '富国沪港深价值精选灵活配置混合A'
$new.Range("D2").Value = '40.15'
$new.Range("E2").Value = '68.33'
$new.Range("F2").Value = '1.75'
$new.Range("G2").Value = '0.7026'
$new.Range("H2").Value = 9

$new.Range("A3").Value = 1
$new.Range("B3").Value = '501025'
$new.Range("C3").Value = '鹏华港股通中证香港银行投资指数（LOF）A'
$new.Range("D3").Value = '9.81'
$new.Range("E3").Value = '94.47'
$new.Range("F3").Value = '5.65'
$new.Range("G3").Value = '0.5543'
$new.Range("H3").Value = 6

$new.Range("A4").Value = 2
$new.Range("B4").Value = '010671'
$new.Range("C4").Value = '景顺长城大中华混合(QDII)美元'
$new.Range("D4").Value = '10.35'
$new.Range("E4").Value = '82.59'
$new.Range("F4").Value = '4.86'
$new.Range("G4").Value = '0.5030'
$new.Range("H4").Value = 8

$new.Range("A5").Value = 3
$new.Range("B5").Value = '262001'
$new.Range("C5").Value = '景顺长城大中华混合(QDII)'
$new.Range("D5").Value = '10.35'
$new.Range("E5").Value = '82.59'
$new.Range("F5").Value = '4.86'
$new.Range("G5").Value = '0.5030'
$new.Range("H5").Value = 8

$new.Range("A6").Value = 4
$new.Range("B6").Value = '010365'
$new.Range("C6").Value = '鹏华港股通中证香港银行投资指数（LOF）C'
$new.Range("D6").Value = '6.07'
$new.Range("E6").Value = '94.47'
$new.Range("F6").Value = '5.65'
$new.Range("G6").Value = '0.3430'
$new.Range("H6").Value = 6

$new.Range("A7").Value = 5
$new.Range("B7").Value = '014746'
$new.Range("C7").Value = '贝莱德港股通远景视野混合A'
$new.Range("D7").Value = '5.05'
$new.Range("E7").Value = '53.79'
$new.Range("F7").Value = '2.29'
$new.Range("G7").Value = '0.1156'
$new.Range("H7").Value = 5

$new.Range("A8").Value = 6
$new.Range("B8").Value = '006809'
$new.Range("C8").Value = '泰康港股通中证香港银行投资指数A'
$new.Range("D8").Value = '1.99'
$new.Range("E8").Value = '94.73'
$new.Range("F8").Value = '5.65'
$new.Range("G8").Value = '0.1124'
$new.Range("H8").Value = 6

$new.Range("A9").Value = 7
$new.Range("B9").Value = '160125'
$new.Range("C9").Value = '南方香港优选股票QDII-LOF'
$new.Range("D9").Value = '2.46'
$new.Range("E9").Value = '91.14'
$new.Range("F9").Value = '3.43'
$new.Range("G9").Value = '0.0844'
$new.Range("H9").Value = 9

$new.Range("A10").Value = 8
$new.Range("B10").Value = '007109'
$new.Range("C10").Value = '南方沪港深核心优势混合'
$new.Range("D10").Value = '1.82'
$new.Range("E10").Value = '87.54'
$new.Range("F10").Value = '2.98'
$new.Range("G10").Value = '0.0542'
$new.Range("H10").Value = 10

$new.Range("A11").Value = 9
$new.Range("B11").Value = '014747'
$new.Range("C11").Value = '贝莱德港股通远景视野混合C'
$new.Range("D11").Value = '2.23'
$new.Range("E11").Value = '53.79'
$new.Range("F11").Value = '2.29'
$new.Range("G11").Value = '0.0511'
$new.Range("H11").Value = 5

$new.Range("A12").Value = 10
$new.Range("B12").Value = '006810'
$new.Range("C12").Value = '泰康港股通中证香港银行投资指数C'
$new.Range("D12").Value = '0.90'
$new.Range("E12").Value = '94.73'
$new.Range("F12").Value = '5.65'
$new.Range("G12").Value = '0.0508'
$new.Range("H12").Value = 6

$new.Range("A13").Value = 11
$new.Range("B13").Value = '011131'
$new.Range("C13").Value = '富国沪港深价值精选灵活配置混合C'
$new.Range("D13").Value = '0.42'
$new.Range("E13").Value = '68.33'
$new.Range("F13").Value = '1.75'
$new.Range("G13").Value = '0.0074'
$new.Range("H13").Value = 9

$new.Range("B2:G13").ClearFormats()

# --- "总计" (totals) sheet ---------------------------------------------------
# Insert a new row 2 for the 2022-Q1 summary line and bump the existing
# rows' running index (column A) down by one.
$tot = $wb.Worksheets.Item("总计")
$tot.Range("A2").EntireRow.Insert()

# The inserted row starts out unformatted; give A2 the same centered/
# bordered style used by the other index cells below it.
$tot.Range("A2:D2").ClearFormats()
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 12
$tot.Range("D2").Value = 3.08

$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
$tot.Range("A7").Value = 5
